{"js": "// \"Added many more features\" \u2014 rewrite the title/meta, refresh the\n// \"What we like\" / \"What we don't like\" bullet lists, and update the\n// closing meta title/description, matching the target diff exactly.\n//\n// We locate each old string with body.search() (exact, case-sensitive)\n// and then replace just that matched sub-range with the new text. Doing\n// the replace on the narrow search-result range (rather than on the\n// whole paragraph) leaves sibling runs \u2014 including the leading empty\n// `<w:r/>` runs used as bookmarks/anchors in the bullet paragraphs, and\n// the bold/italic runs near the end \u2014 untouched.\n\nasync function replaceAll(context, oldText, newText) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Title (Heading1 at the top, and the bold \"meta title\" line near the end)\nawait replaceAll(\n  context,\n  \"Play Mayan Riches Online Slot for Free | Review\",\n  \"Play Mayan Riches Free - Exciting Mayan-themed Slot Game\"\n);\n\n// \"What we like\" bullets\nawait replaceAll(\n  context,\n  \"Stackable Wild symbols during free spins\",\n  \"Wide range of betting options\"\n);\nawait replaceAll(\n  context,\n  \"Exceptional graphics and vibrant colors\",\n  \"Stackable wild symbols during free spins\"\n);\nawait replaceAll(\n  context,\n  \"Multi-platform accessibility\",\n  \"Exceptional graphics and vivid colors\"\n);\nawait replaceAll(\n  context,\n  \"Plenty of free spins available\",\n  \"Mesmerizing soundtrack\"\n);\n\n// \"What we don't like\" bullets\nawait replaceAll(\n  context,\n  \"Mayan theme not original\",\n  \"Mayan theme may not appeal to everyone\"\n);\nawait replaceAll(\n  context,\n  \"No bonus game feature\",\n  \"Limited variety of bonus features\"\n);\n\n// Meta description (italic run near the end)\nawait replaceAll(\n  context,\n  \"Read our review of Mayan Riches, a classic online slot with 40 paylines and a tribal Mayan theme. Play for free and enjoy the exceptional graphics and sound design.\",\n  \"Play Mayan Riches for free and immerse yourself in the ancient Maya culture.\"\n);\n", "ps1": "# Applies the \"Added many more features\" edit to the Mayan Riches review.\n# Uses Find/Replace (wdReplaceAll) over the whole document content so that\n# each paragraph's existing run formatting (bold/italic/etc.) is preserved.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All([string]$findText, [string]$replaceText) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $findText\n  $rng.Find.Replacement.Text = $replaceText\n  $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Title (appears both as the Heading1 at the top and as the bold \"meta title\" near the end)\nReplace-All \"Play Mayan Riches Online Slot for Free | Review\" \"Play Mayan Riches Free - Exciting Mayan-themed Slot Game\"\n\n# \"What we like\" bullets\nReplace-All \"Stackable Wild symbols during free spins\" \"Wide range of betting options\"\nReplace-All \"Exceptional graphics and vibrant colors\" \"Stackable wild symbols during free spins\"\nReplace-All \"Multi-platform accessibility\" \"Exceptional graphics and vivid colors\"\nReplace-All \"Plenty of free spins available\" \"Mesmerizing soundtrack\"\n\n# \"What we don't like\" bullets\nReplace-All \"Mayan theme not original\" \"Mayan theme may not appeal to everyone\"\nReplace-All \"No bonus game feature\" \"Limited variety of bonus features\"\n\n# Meta description (italic run near the end)\nReplace-All \"Read our review of Mayan Riches, a classic online slot with 40 paylines and a tribal Mayan theme. Play for free and enjoy the exceptional graphics and sound design.\" \"Play Mayan Riches for free and immerse yourself in the ancient Maya culture.\"\n"}
